# Estadisticos Segundo Parcial 23 Mayo
$wb = $excel.ActiveWorkbook

# --- "Estadisticos 2P" sheet: update stats row (student passed/failed numbers) ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 17
$ws2P.Range("E2").Value = 17
$ws2P.Range("F2").Value = 22
$ws2P.Range("G2").Value = 56.41
$ws2P.Range("H2").Value = 9.5

# --- "Estadisticos Final" sheet: refresh the average (Promedio) ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("H2").Value = 9.1

# --- "Rescatables" sheet: add the two rescatable students ---
$wsResc = $wb.Worksheets.Item("Rescatables")

$wsResc.Cells.Item(2, 1).Value = 21330051920053
$wsResc.Cells.Item(3, 1).Value = 21330051920242

$wsResc.Cells.Item(2, 2).Value = "OREA"
$wsResc.Cells.Item(3, 2).Value = "PALMA"

$wsResc.Cells.Item(2, 3).Value = "MARTINEZ"
$wsResc.Cells.Item(3, 3).Value = "RANGEL"

$wsResc.Cells.Item(2, 4).Value = "JOSE MANUEL"
$wsResc.Cells.Item(3, 4).Value = "ROBERTO"

$wsResc.Cells.Item(2, 5).Value = "TEMAS DE ADMINISTRACIÓN"
$wsResc.Cells.Item(3, 5).Value = "TEMAS DE ADMINISTRACIÓN"

$wsResc.Cells.Item(2, 6).Value = "6ARHM"
$wsResc.Cells.Item(3, 6).Value = "6ARHM"

$wsResc.Cells.Item(2, 7).Value = 4
$wsResc.Cells.Item(3, 7).Value = 4
